$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh with new TPM-derived values, target cluster -> MuSCs ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vip"
$ws.Range("C2").Value = "Sctr"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6855646666666667
$ws.Range("H2").Value = 2.056694
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.015638
$ws.Range("N2").Value = 0.046914
$ws.Range("O2").Value = 0.05572514913016313
$ws.Range("P2").Value = 0.05572514913016313
$ws.Range("Q2").Value = 0.01072086025733333
$ws.Range("R2").Value = 0.09648774231600001
$ws.Range("S2").Value = 0.05572514913016313
$ws.Range("T2").Value = 0.05572514913016313

# --- Row 3 (new): same sender/ligand/receptor, target cluster -> Resolving-Mac ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vip"
$ws.Range("C3").Value = "Sctr"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6855646666666667
$ws.Range("H3").Value = 2.056694
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2649893333333334
$ws.Range("N3").Value = 0.794968
$ws.Range("O3").Value = 0.9442748508698369
$ws.Range("P3").Value = 0.9442748508698369
$ws.Range("Q3").Value = 0.1816673239768889
$ws.Range("R3").Value = 1.635005915792
$ws.Range("S3").Value = 0.9442748508698369
$ws.Range("T3").Value = 0.9442748508698369
